$d = $word.ActiveDocument

$replacements = @(
    @("244÷6=", "192÷7="),
    @("392÷6=", "263÷4="),
    @("595÷9=", "951÷9="),
    @("280÷7=", "519÷8="),
    @("891÷8=", "486÷4="),
    @("391÷8=", "963÷3="),
    @("349÷5=", "780÷5="),
    @("301÷6=", "571÷6="),
    @("459÷4=", "521÷4="),
    @("369÷7=", "860÷5="),
    @("326÷4=", "807÷2="),
    @("113÷6=", "784÷8="),
    @("503÷4=", "836÷3="),
    @("543÷7=", "602÷6="),
    @("514÷5=", "164÷4="),
    @("779÷2=", "819÷4="),
    @("542÷3=", "618÷3="),
    @("905÷4=", "748÷8="),
    @("729÷5=", "481÷6="),
    @("413÷8=", "502÷2="),
    @("398÷2=", "378÷5="),
    @("594÷2=", "437÷3="),
    @("352÷4=", "932÷7="),
    @("833÷8=", "837÷4="),
    @("808÷4=", "529÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
